$d = $word.ActiveDocument

# The paragraph currently reads "First section second column column" (a
# duplicated trailing "column"). Trim it down to "First section second column".
$find = $d.Content
$find.Find.Execute("First section second column column", $true, $false, $false, $false, $false, $true, 1, $false, "First section second column", 2)

# Re-locate the (now corrected) run so we can split it into three runs, as in
# the authored edit: "First section second " + "c" + "olumn".
$target = $d.Content
$target.Find.Execute("First section second column", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $target.Start

# "First section second " is 21 characters (offsets $start .. $start+21),
# followed by a lone "c" (1 character), followed by "olumn".
$splitPoint1 = $start + 21
$splitPoint2 = $start + 22

# Force a run break around the single "c" character by nudging a character
# formatting property on and back off; Word splits runs at range boundaries
# whenever a run-level property is applied to a sub-range.
$cRange = $d.Range($splitPoint1, $splitPoint2)
$cRange.Bold = 1
$cRange.Bold = 0
